# Updated cryptos list on Thu Feb 29 15:21:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.917.10"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").Value = "3.472.85"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +17.73%  "

$ws.Range("D7").Value = "3.466.45"
$ws.Range("E7").Value = "  +2.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.699"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.131"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +31.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.11%  "

$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").Value = "4.011.51"
$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("E15").Value = "  +4.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.16"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "3.448.36"

$ws.Range("D18").Value = "62.809.75"
$ws.Range("E18").Value = "  +3.18%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000144"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +29.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "312.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.92%  "

$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.62%  "

$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.177"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.64%  "

$ws.Range("E32").Value = "  +2.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "43.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0495"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.54%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("E43").Value = "  +3.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.290"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("D50").Value = "3.810.75"
$ws.Range("E50").Value = "  +3.44%  "

$ws.Range("D51").Value = "2.183.28"
$ws.Range("E51").Value = "  -0.02%  "
